$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 20) of forecast data, matching the style of the
# previous row (row 19) for the date cell in column A.
$ws.Range("A19:E19").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)

$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 1.049317648994741
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 0.3243937446859801
